## Updated SQL exports + saving progress
## Fills in the "facility_group_name" (and new "comment") classification for the
## facility_groups sheet: Glencore's Sudbury-area smelters/mines get grouped
## under "Sudbury INO" / "Integrated Nickel Operations", ArcelorMittal's three
## Contrecoeur plants get grouped under "Contrecoeur", and ArcelorMittal's
## Fire Lake / Mont-Wright / Port-Cartier trio gets grouped under "AMMC" with
## a shared sourcing comment. Also documents Glencore Horne's CCR? grouping
## with its own comment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("facility_groups")

# --- 1. Rows 36:39 were a blank gap between the Vale block (rows 29-35) and
#     the Glencore block (old rows 40-49). Delete the gap so the Glencore
#     block slides up to rows 36-45, matching the tidied-up layout. ---
$ws.Rows("36:39").Delete()

# --- 2. Fill in the facility_group_name (column D) for the existing Glencore
#     rows (now at rows 36-45), in the order the author actually typed them:
#     Integrated Nickel Operations first, then Sudbury INO group, then the
#     CCR group, then the new rows/columns added below. ---

# Raglan (row 43) and Sudbury (row 45) -> "Integrated Nickel Operations"
$ws.Range("D43").Value = "Integrated Nickel Operations"
$ws.Range("D45").Value = "Integrated Nickel Operations"

# Fraser (row 38), Nickel Rim South (row 42), Strathcona (row 44) -> "Sudbury INO"
$ws.Range("D38").Value = "Sudbury INO"
$ws.Range("D42").Value = "Sudbury INO"
$ws.Range("D44").Value = "Sudbury INO"

# CCR (row 37) -> "CCR" (re-uses the already-existing "CCR" shared string)
$ws.Range("D37").Value = "CCR"

# --- 3. Append the three new ArcelorMittal Long Products Contrecoeur rows
#     (46-48), each grouped under "Contrecoeur". ---
$ws.Range("A46").Value = "ArcelorMittal Long Products Canada"
$ws.Range("B46").Value = "Contrecœur East"
$ws.Range("C46").Value = "QC-MAIN-1eee4ace"
$ws.Range("D46").Value = "Contrecoeur"

$ws.Range("A47").Value = "ArcelorMittal Long Products Canada"
$ws.Range("B47").Value = "Contrecœur West"
$ws.Range("C47").Value = "QC-MAIN-844dcc47"
$ws.Range("D47").Value = "Contrecoeur"

$ws.Range("A48").Value = "ArcelorMittal Long Products Canada"
$ws.Range("B48").Value = "Contrecœur-Feruni"
$ws.Range("C48").Value = "QC-MAIN-30ff61c3"
$ws.Range("D48").Value = "Contrecoeur"

# --- 4. Append the three new ArcelorMittal AMMC rows (49-51): Fire Lake,
#     Mont-Wright and Port-Cartier, all grouped under "AMMC". ---
$ws.Range("A49").Value = "ArcelorMittal"
$ws.Range("B49").Value = "Fire Lake"
$ws.Range("C49").Value = "QC-MAIN-084bd95c"
$ws.Range("D49").Value = "AMMC"

$ws.Range("A50").Value = "ArcelorMittal"
$ws.Range("B50").Value = "Mont-Wright"
$ws.Range("C50").Value = "QC-MAIN-33c09b8b"
$ws.Range("D50").Value = "AMMC"

$ws.Range("A51").Value = "ArcelorMittal Mining Canada"
$ws.Range("B51").Value = "Port-Cartier"
$ws.Range("C51").Value = "QC-MAIN-e25eed27"
$ws.Range("D51").Value = "AMMC"

# --- 5. Add the new "comment" column header. ---
$ws.Range("F1").Value = "comment"

# --- 6. Sourcing comment (quoted verbatim from ArcelorMittal) shared by the
#     three AMMC rows, entered in column F, italicised. ---
$fireLakeComment = @"
"The raw ore from Fire Lake is sent to our Mont-Wright mining complex where all our ore is crushed, ground and concentrated. Our concentrate is then transported by rail to Port-Cartier, where it will either be sent to our plant to be made into iron oxide pellets or directly to our port to be shipped to international markets.
The combined production of our Mont-Wright and Fire Lake mines represents over 30% of the ArcelorMittal Group’s global iron ore supply."
"@

$ws.Range("F49").Value = $fireLakeComment
$ws.Range("F49").Font.Italic = $true
$ws.Range("F50").Value = $fireLakeComment
$ws.Range("F50").Font.Italic = $true
$ws.Range("F51").Value = $fireLakeComment
$ws.Range("F51").Font.Italic = $true

# --- 7. Horne (row 40): tentative "CCR?" grouping plus its own sourcing
#     comment, also italicised. ---
$ws.Range("D40").Value = "CCR?"

$horneComment = @"
"At Horne Smelter, we produce 99.1% pure copper anodes. From their arrival at our Rouyn-Noranda smelter, whether by train or truck, to casting in our anode furnaces, the concentrates and recycled products go through various processing stages before the material can then be processed into cathodes at the CCR Refinery in Montréal-Est."
"@

$ws.Range("F40").Value = $horneComment
$ws.Range("F40").Font.Italic = $true

# --- 8. Leave the view parked on the newly-edited block, like the author did
#     while reviewing their work. ---
$ws.Range("A29:D35").Select()
